$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need an explicit Text format
# so Excel stores them as text (matching the source inline-string cells)
# instead of silently coercing them into numbers.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.238.28"
$ws.Range("E2").Value = "  +5.49%  "
$ws.Range("D3").Value = "1.883.17"
$ws.Range("E3").Value = "  +3.87%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "282.11"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("D7").Value = "0.5314"
$ws.Range("E7").Value = "  +4.39%  "
$ws.Range("D8").Value = "0.3540"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "0.07059"
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").Value = "20.46"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "0.8235"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "0.07819"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.897.13"
$ws.Range("E13").Value = "  +4.73%  "
$ws.Range("D14").Value = "91.05"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").Value = "5.214"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "0.9977"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "14.64"
$ws.Range("E17").Value = "  +5.40%  "
$ws.Range("D19").Value = "0.9991"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.285.14"
$ws.Range("E20").Value = "  +5.43%  "
$ws.Range("D21").Value = "2.125.47"
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "4.780"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "6.257"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").Value = "2.407"
$ws.Range("E25").Value = "  +8.72%  "
$ws.Range("D26").Value = "147.18"
$ws.Range("E26").Value = "  +3.82%  "
$ws.Range("D27").Value = "17.65"
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("D28").Value = "1.680"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").Value = "114.42"
$ws.Range("E29").Value = "  +5.08%  "
$ws.Range("D30").Value = "4.427"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").Value = "4.400"
$ws.Range("E31").Value = "  +4.29%  "
$ws.Range("D32").Value = "0.08966"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "0.04953"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "1.188"
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("D35").Value = "0.7506"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").Value = "2.905"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "3.319"
$ws.Range("E37").Value = "  +9.11%  "
$ws.Range("D38").Value = "2.438"
$ws.Range("E38").Value = "  +6.17%  "
$ws.Range("D39").Value = "0.5341"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").Value = "0.01889"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "0.9748"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("D42").Value = "117.07"
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("D43").Value = "6.341"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").Value = "8.266"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "0.4629"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").Value = "0.9980"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "0.1374"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "9.470"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "36.83"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "0.05959"
$ws.Range("E51").Value = "  +2.72%  "

Write-Output "Updated cryptos list"
